$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Strip leading whitespace from the LCSC Part # column (D2:D7), keeping
# trailing spaces intact.
$ws.Range("D2").Value = "C15849     "
$ws.Range("D3").Value = "C14663     "
$ws.Range("D4").Value = "C190321    "
$ws.Range("D5").Value = "C23162     "
$ws.Range("D6").Value = "C115357    "
$ws.Range("D7").Value = "C7377      "

# Update the saved selection / active cell for the sheet view.
$ws.Range("C16").Select()
